$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de report sheets (regenerated report run).

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 03:28:11"
$wsZhCn.Range("H2").Value = "2016-03-25 03:28:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 03:28:16"
$wsDeDe.Range("H2").Value = "2016-03-25 03:28:44"
